$wb = $excel.ActiveWorkbook

# Insert three new sheets right after "Select Input" and before "Radio Buttons Demo":
#   FormSubmit, FormSubmitByParameter, FormSubmitByDataProvider
$selectInput = $wb.Worksheets.Item("Select Input")

$formSubmit = $wb.Worksheets.Add($null, $selectInput)
$formSubmit.Name = "FormSubmit"
$formSubmit.Range("A1").Value = "First Name"
$formSubmit.Range("B1").Value = "Last Name"
$formSubmit.Range("C1").Value = "UserName"
$formSubmit.Range("D1").Value = "City"
$formSubmit.Range("E1").Value = "State"
$formSubmit.Range("F1").Value = "Zip"
$formSubmit.Range("A2").Value = "Sneha"
$formSubmit.Range("B2").Value = "George"
$formSubmit.Range("C2").Value = "Sneha95"
$formSubmit.Range("D2").Value = "Kottayam"
$formSubmit.Range("E2").Value = "Kerala"
$formSubmit.Range("F2").Value = 698473
$formSubmit.Range("G7").Select() | Out-Null

$formSubmitByParameter = $wb.Worksheets.Add($null, $formSubmit)
$formSubmitByParameter.Name = "FormSubmitByParameter"
$formSubmitByParameter.Range("A1").Value = "Expected String"
$formSubmitByParameter.Range("A2").Value = "Form has been submitted successfully!"
$formSubmitByParameter.Range("A2").Select() | Out-Null

$formSubmitByDataProvider = $wb.Worksheets.Add($null, $formSubmitByParameter)
$formSubmitByDataProvider.Name = "FormSubmitByDataProvider"
$formSubmitByDataProvider.Range("A1").Value = "Expected String"
$formSubmitByDataProvider.Range("A2").Value = "Form has been submitted successfully!"
$formSubmitByDataProvider.Range("A2").Select() | Out-Null

# The newly created, currently-active sheet (FormSubmitByDataProvider) becomes the
# selected tab, matching the target workbook's activeTab/tabSelected state.
